# Append 5 new tracker rows (208-212) to the mortality data sheet,
# matching the "Data updated 06 May 2024" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns A and B hold dd.mm.yyyy date strings stored as plain TEXT
# (shared strings) in the source file, not as real Excel dates. Force
# text interpretation for the whole block first so "02.05.2024" etc.
# aren't auto-converted into date serials, then restore General so the
# cells keep the same (default) style as the rest of the sheet.
$dateTextRange = $ws.Range("A208:B212")
$dateTextRange.NumberFormat = "@"

# --- Row 208 : 02.05.2024 (report_date 02.05.2024) ---
$ws.Range("A208").Value = "02.05.2024"
$ws.Range("B208").Value = "02.05.2024"
$ws.Range("C208").Value = 34596
$ws.Range("D208").Value = 14500
$ws.Range("E208").Value = 8400
$ws.Range("F208").Value = 77816
$ws.Range("I208").Value = 8000
$ws.Range("J208").Value = 492
$ws.Range("K208").Value = 124
$ws.Range("L208").Value = 4800
$ws.Range("M208").Value = "https://web.archive.org/web/20240502164235/https://www.aljazeera.com/news/longform/2023/10/9/israel-hamas-war-in-maps-and-charts-live-tracker"

# --- Row 209 : 03.05.2024 (report_date 02.05.2024) ---
$ws.Range("A209").Value = "03.05.2024"
$ws.Range("B209").Value = "02.05.2024"
$ws.Range("C209").Value = 34622
$ws.Range("D209").Value = 14500
$ws.Range("E209").Value = 8400
$ws.Range("F209").Value = 77867
$ws.Range("I209").Value = 8000
$ws.Range("J209").Value = 492
$ws.Range("K209").Value = 124
$ws.Range("L209").Value = 4800
$ws.Range("M209").Value = "https://web.archive.org/web/20240503142834/https://www.aljazeera.com/news/longform/2023/10/9/israel-hamas-war-in-maps-and-charts-live-tracker"

# --- Row 210 : 04.05.2024 (report_date 02.05.2024) ---
$ws.Range("A210").Value = "04.05.2024"
$ws.Range("B210").Value = "02.05.2024"
$ws.Range("C210").Value = 34622
$ws.Range("D210").Value = 14500
$ws.Range("E210").Value = 8400
$ws.Range("F210").Value = 77867
$ws.Range("I210").Value = 8000
$ws.Range("J210").Value = 492
$ws.Range("K210").Value = 124
$ws.Range("L210").Value = 4800
$ws.Range("M210").Value = "https://web.archive.org/web/20240504232915/https://www.aljazeera.com/news/longform/2023/10/9/israel-hamas-war-in-maps-and-charts-live-tracker"

# --- Row 211 : 05.04.2024 (reused date string) (report_date 02.05.2024) ---
$ws.Range("A211").Value = "05.04.2024"
$ws.Range("B211").Value = "02.05.2024"
$ws.Range("C211").Value = 34622
$ws.Range("D211").Value = 14500
$ws.Range("E211").Value = 8400
$ws.Range("F211").Value = 77867
$ws.Range("I211").Value = 8000
$ws.Range("J211").Value = 492
$ws.Range("K211").Value = 124
$ws.Range("L211").Value = 4800
$ws.Range("M211").Value = "https://web.archive.org/web/20240505201358/https://www.aljazeera.com/news/longform/2023/10/9/israel-hamas-war-in-maps-and-charts-live-tracker"

# --- Row 212 : 06.05.2024 (report_date 06.05.2024) ---
$ws.Range("A212").Value = "06.05.2024"
$ws.Range("B212").Value = "06.05.2024"
$ws.Range("C212").Value = 34735
$ws.Range("D212").Value = 14500
$ws.Range("E212").Value = 8400
$ws.Range("F212").Value = 78108
$ws.Range("I212").Value = 8000
$ws.Range("J212").Value = 497
$ws.Range("K212").Value = 124
$ws.Range("L212").Value = 4800
$ws.Range("M212").Value = "https://web.archive.org/web/20240506141604/https://www.aljazeera.com/news/longform/2023/10/9/israel-hamas-war-in-maps-and-charts-live-tracker"

# Restore the default General format on A208:B212 now that the values
# are committed as text, so the cells keep the sheet's default style.
$dateTextRange.NumberFormat = "General"

# Row 200-207 alternate an explicit wrap-text style on column C; the new
# block of rows keeps that pattern going on the first new row (C208).
$ws.Range("C208").WrapText = $true

# Match the author's final cursor position (bottom pane, cell M211).
$null = $ws.Range("M211").Select()
